# Updated cryptos list values (price + volume change columns) to match
# the latest scrape. Column D (Price) entries that look numeric get a
# leading apostrophe so Excel keeps them as literal text (matching the
# original inlineStr cells) instead of silently coercing to a Double and
# dropping significant trailing zeros (e.g. "1.150" -> 1.15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.246.82"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "1.803.28"
$ws.Range("E3").Value = "  +2.75%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'336.86"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4637"
$ws.Range("E7").Value = "  +20.71%  "

$ws.Range("D8").Value = "'0.3758"
$ws.Range("E8").Value = "  +10.19%  "

$ws.Range("D9").Value = "'45.17"
$ws.Range("E9").Value = "  -2.36%  "

$ws.Range("D10").Value = "'0.07637"
$ws.Range("E10").Value = "  +5.43%  "

$ws.Range("D11").Value = "'1.150"
$ws.Range("E11").Value = "  +2.40%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "'22.36"
$ws.Range("E13").Value = "  -1.18%  "

$ws.Range("D14").Value = "'6.330"
$ws.Range("E14").Value = "  +2.44%  "

$ws.Range("D15").Value = "'7.460"
$ws.Range("E15").Value = "  +4.18%  "

$ws.Range("D16").Value = "1.808.40"
$ws.Range("E16").Value = "  +3.20%  "

$ws.Range("D17").Value = "'0.00001096"
$ws.Range("E17").Value = "  +3.16%  "

$ws.Range("D18").Value = "'0.06725"
$ws.Range("E18").Value = "  +1.75%  "

$ws.Range("D19").Value = "'81.84"
$ws.Range("E19").Value = "  +3.29%  "

$ws.Range("D20").Value = "'0.9996"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").Value = "'17.43"
$ws.Range("E21").Value = "  +3.83%  "

$ws.Range("D22").Value = "'6.423"
$ws.Range("E22").Value = "  +2.94%  "

$ws.Range("D23").Value = "28.240.37"
$ws.Range("E23").Value = "  +1.05%  "

$ws.Range("D24").Value = "'11.85"
$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("E25").Value = "  +1.13%  "

$ws.Range("D26").Value = "'20.78"
$ws.Range("E26").Value = "  +4.44%  "

$ws.Range("D27").Value = "'154.02"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").Value = "'2.375"
$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("D29").Value = "2.009.92"
$ws.Range("E29").Value = "  +2.93%  "

$ws.Range("D30").Value = "'133.42"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("E31").Value = "  -1.93%  "

$ws.Range("D32").Value = "'4.033"
$ws.Range("E32").Value = "  +0.41%  "

$ws.Range("D33").Value = "'0.09594"
$ws.Range("E33").Value = "  +8.64%  "

$ws.Range("D34").Value = "'5.856"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").Value = "'0.2232"
$ws.Range("E35").Value = "  +5.80%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06376"
$ws.Range("E36").Value = "  +2.88%  "

$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'12.11"
$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.271"
$ws.Range("E38").Value = "  +2.01%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02349"
$ws.Range("E39").Value = "  +2.30%  "

$ws.Range("D40").Value = "'0.6636"
$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("D41").Value = "'1.236"
$ws.Range("E41").Value = "  +1.74%  "

$ws.Range("D42").Value = "'1.497"
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").Value = "'8.260"
$ws.Range("E43").Value = "  +3.16%  "

$ws.Range("D44").Value = "'14.19"
$ws.Range("E44").Value = "  +3.04%  "

$ws.Range("D45").Value = "'0.9995"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("D47").Value = "'3.830"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").Value = "'129.76"
$ws.Range("E48").Value = "  +2.03%  "

$ws.Range("D49").Value = "'2.040"
$ws.Range("E49").Value = "  +1.21%  "

$ws.Range("D50").Value = "'0.07162"
$ws.Range("E50").Value = "  +2.56%  "

$ws.Range("D51").Value = "'1.175"
$ws.Range("E51").Value = "  -0.22%  "
